$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, pushing the existing data (rows 102-147) down to 103-148
$ws.Rows(102).Insert()

# Populate the newly inserted row 102 with a new weekly price record
$ws.Range("A102").Value = 4
$ws.Range("B102").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C102").Value = "Los Lagos"
$ws.Range("D102").Value = 44572
$ws.Range("E102").Value = 10
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100108
$ws.Range("H102").Value = "Tropicales y subtropicales"
$ws.Range("I102").Value = 100108002
$ws.Range("J102").Value = "Mango"
$ws.Range("K102").Value = "Sin especificar"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 200
$ws.Range("N102").Value = 8000
$ws.Range("O102").Value = 8500
$ws.Range("P102").Value = 8250
$ws.Range("Q102").Value = "$/bandeja 4 kilos"
$ws.Range("R102").Value = "Perú"
$ws.Range("S102").Value = 2062
$ws.Range("T102").Value = 4
